$wb = $excel.ActiveWorkbook

# ---- Create sheet 03-01-2024 (copy of 02-29-2024) ----
$src1129 = $wb.Worksheets.Item("02-29-2024")
$src1129.Copy($null, $src1129)
$sheet0301 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet0301.Name = "03-01-2024"

# Update numeric predicted-difference values for 03-01-2024
$sheet0301.Cells.Item(2, 2).Value = 2.55735312101535
$sheet0301.Cells.Item(2, 3).Value = 1.428045749664307
$sheet0301.Cells.Item(3, 2).Value = 3.552435299359359
$sheet0301.Cells.Item(3, 3).Value = 2.704452753067017
$sheet0301.Cells.Item(4, 2).Value = 0.1573702929926828
$sheet0301.Cells.Item(4, 3).Value = -0.6290268301963806
$sheet0301.Cells.Item(5, 2).Value = -1.685259414014917
$sheet0301.Cells.Item(5, 3).Value = -2.005419492721558
$sheet0301.Cells.Item(6, 2).Value = 0.4935147355608027
$sheet0301.Cells.Item(6, 3).Value = -0.002308783587068319
$sheet0301.Cells.Item(7, 2).Value = 2.266349808910419
$sheet0301.Cells.Item(7, 3).Value = 0.6656510233879089
$sheet0301.Cells.Item(8, 2).Value = 1.896363546491922
$sheet0301.Cells.Item(8, 3).Value = 2.818512916564941
$sheet0301.Cells.Item(9, 2).Value = -1.623009936314896
$sheet0301.Cells.Item(9, 3).Value = -2.209126472473145
$sheet0301.Cells.Item(10, 2).Value = 0.6786593885306607
$sheet0301.Cells.Item(10, 3).Value = 1.1358642578125
$sheet0301.Cells.Item(11, 2).Value = 6.059013859872003
$sheet0301.Cells.Item(11, 3).Value = 2.914090871810913
$sheet0301.Cells.Item(12, 2).Value = 5.420697877993796
$sheet0301.Cells.Item(12, 3).Value = 4.099342823028564
$sheet0301.Cells.Item(13, 2).Value = 1.045554642035351
$sheet0301.Cells.Item(13, 3).Value = 0.8425951600074768
$sheet0301.Cells.Item(14, 2).Value = 5.544195999989706
$sheet0301.Cells.Item(14, 3).Value = 4.250433444976807
$sheet0301.Cells.Item(15, 2).Value = 4.208166318466964
$sheet0301.Cells.Item(15, 3).Value = 7.006680965423584
$sheet0301.Cells.Item(16, 2).Value = -0.01151784714039916
$sheet0301.Cells.Item(16, 3).Value = 0.9264268279075623
$sheet0301.Cells.Item(17, 2).Value = 3.100889504997216
$sheet0301.Cells.Item(17, 3).Value = 1.376412630081177
$sheet0301.Cells.Item(18, 2).Value = 0.7862631708123846
$sheet0301.Cells.Item(18, 3).Value = 1.009001135826111
$sheet0301.Cells.Item(19, 2).Value = 1.440953942670495
$sheet0301.Cells.Item(19, 3).Value = 3.303689002990723
$sheet0301.Cells.Item(20, 2).Value = -1.239776003188227
$sheet0301.Cells.Item(20, 3).Value = -1.770082592964172
$sheet0301.Cells.Item(21, 2).Value = -0.5649538828080178
$sheet0301.Cells.Item(21, 3).Value = -2.049990892410278
$sheet0301.Cells.Item(22, 2).Value = 3.85324205731877
$sheet0301.Cells.Item(22, 3).Value = 3.715561151504517
$sheet0301.Cells.Item(23, 2).Value = -2.106638662432982
$sheet0301.Cells.Item(23, 3).Value = -1.068375587463379
$sheet0301.Cells.Item(24, 2).Value = 0.7605801528544447
$sheet0301.Cells.Item(24, 3).Value = 0.6613050103187561
$sheet0301.Cells.Item(25, 2).Value = -1.40710375698842
$sheet0301.Cells.Item(25, 3).Value = -5.379960536956787
$sheet0301.Cells.Item(26, 2).Value = 6.077012849577852
$sheet0301.Cells.Item(26, 3).Value = 5.243505954742432
$sheet0301.Cells.Item(27, 2).Value = 2.517180178338744
$sheet0301.Cells.Item(27, 3).Value = 1.322910666465759
$sheet0301.Cells.Item(28, 2).Value = 0.7425910700587028
$sheet0301.Cells.Item(28, 3).Value = 1.041263937950134
$sheet0301.Cells.Item(29, 2).Value = -11.18040598727329
$sheet0301.Cells.Item(29, 3).Value = -7.766969203948975
$sheet0301.Cells.Item(30, 2).Value = 3.998123619421742
$sheet0301.Cells.Item(30, 3).Value = 3.695779085159302
$sheet0301.Cells.Item(31, 2).Value = 6.140934624196198
$sheet0301.Cells.Item(31, 3).Value = 5.420891284942627
$sheet0301.Cells.Item(32, 2).Value = -0.6806780331245896
$sheet0301.Cells.Item(32, 3).Value = 0.6375436782836914
$sheet0301.Cells.Item(33, 2).Value = -0.102497547776732
$sheet0301.Cells.Item(33, 3).Value = -2.372321844100952
$sheet0301.Cells.Item(34, 2).Value = -0.761712023890305
$sheet0301.Cells.Item(34, 3).Value = 0.3949662446975708
$sheet0301.Cells.Item(35, 2).Value = -1.683884078036758
$sheet0301.Cells.Item(35, 3).Value = -3.230232238769531

# ---- Create sheet 03-02-2024 (copy of 02-29-2024, with one extra pollster row) ----
$src1129b = $wb.Worksheets.Item("02-29-2024")
$src1129b.Copy($null, $sheet0301)
$sheet0302 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet0302.Name = "03-02-2024"

# Insert new row for pollster_Forbes/HarrisX at row 10 (shifts rows 10-35 down to 11-36)
$sheet0302.Rows.Item(10).Insert()
$sheet0302.Range("A9").Copy()
$sheet0302.Range("A10").PasteSpecial(-4122)
$sheet0302.Range("A10").Value = "pollster_Forbes/HarrisX"

# Update numeric predicted-difference values for 03-02-2024
$sheet0302.Cells.Item(2, 2).Value = 2.596673397619933
$sheet0302.Cells.Item(2, 3).Value = 4.145450592041016
$sheet0302.Cells.Item(3, 2).Value = 3.591454679856594
$sheet0302.Cells.Item(3, 3).Value = 5.368276119232178
$sheet0302.Cells.Item(4, 2).Value = 0.1722176861934672
$sheet0302.Cells.Item(4, 3).Value = 1.798820734024048
$sheet0302.Cells.Item(5, 2).Value = -1.66078334537639
$sheet0302.Cells.Item(5, 3).Value = 0.1950432509183884
$sheet0302.Cells.Item(6, 2).Value = 0.5466924103169304
$sheet0302.Cells.Item(6, 3).Value = 1.431107759475708
$sheet0302.Cells.Item(7, 2).Value = 2.333753722207729
$sheet0302.Cells.Item(7, 3).Value = 3.038323640823364
$sheet0302.Cells.Item(8, 2).Value = 1.944189153066569
$sheet0302.Cells.Item(8, 3).Value = 5.414598941802979
$sheet0302.Cells.Item(9, 2).Value = -1.533540308473303
$sheet0302.Cells.Item(9, 3).Value = 0.214027464389801
$sheet0302.Cells.Item(10, 2).Value = 4.020874871053748
$sheet0302.Cells.Item(10, 3).Value = 3.995976448059082
$sheet0302.Cells.Item(11, 2).Value = 0.725401769118212
$sheet0302.Cells.Item(11, 3).Value = 3.694656133651733
$sheet0302.Cells.Item(12, 2).Value = 6.067843330924676
$sheet0302.Cells.Item(12, 3).Value = 4.926626205444336
$sheet0302.Cells.Item(13, 2).Value = 5.491662395228951
$sheet0302.Cells.Item(13, 3).Value = 6.351988315582275
$sheet0302.Cells.Item(14, 2).Value = 1.078034547131181
$sheet0302.Cells.Item(14, 3).Value = 3.249961137771606
$sheet0302.Cells.Item(15, 2).Value = 5.643896095483688
$sheet0302.Cells.Item(15, 3).Value = 5.89842414855957
$sheet0302.Cells.Item(16, 2).Value = 4.256714075698859
$sheet0302.Cells.Item(16, 3).Value = 9.394993782043457
$sheet0302.Cells.Item(17, 2).Value = 0.05418098821487671
$sheet0302.Cells.Item(17, 3).Value = 3.618716716766357
$sheet0302.Cells.Item(18, 2).Value = 3.17346596566497
$sheet0302.Cells.Item(18, 3).Value = 4.074489116668701
$sheet0302.Cells.Item(19, 2).Value = 0.8284647555971447
$sheet0302.Cells.Item(19, 3).Value = 3.218261241912842
$sheet0302.Cells.Item(20, 2).Value = 1.488448563526754
$sheet0302.Cells.Item(20, 3).Value = 5.681851387023926
$sheet0302.Cells.Item(21, 2).Value = -1.203339314279646
$sheet0302.Cells.Item(21, 3).Value = 1.047272682189941
$sheet0302.Cells.Item(22, 2).Value = 0.2462766401719474
$sheet0302.Cells.Item(22, 3).Value = 3.998223304748535
$sheet0302.Cells.Item(23, 2).Value = 3.910666249720252
$sheet0302.Cells.Item(23, 3).Value = 5.489115238189697
$sheet0302.Cells.Item(24, 2).Value = -1.985574943312053
$sheet0302.Cells.Item(24, 3).Value = 1.7064288854599
$sheet0302.Cells.Item(25, 2).Value = 0.8735196770749796
$sheet0302.Cells.Item(25, 3).Value = 3.102028369903564
$sheet0302.Cells.Item(26, 2).Value = -1.342076151411259
$sheet0302.Cells.Item(26, 3).Value = -3.36479663848877
$sheet0302.Cells.Item(27, 2).Value = 6.134008843672303
$sheet0302.Cells.Item(27, 3).Value = 7.477080821990967
$sheet0302.Cells.Item(28, 2).Value = 2.569338688860061
$sheet0302.Cells.Item(28, 3).Value = 3.978230237960815
$sheet0302.Cells.Item(29, 2).Value = 0.7932451000428746
$sheet0302.Cells.Item(29, 3).Value = 3.724176645278931
$sheet0302.Cells.Item(30, 2).Value = -11.06385570976372
$sheet0302.Cells.Item(30, 3).Value = -5.112802982330322
$sheet0302.Cells.Item(31, 2).Value = 4.049116157526757
$sheet0302.Cells.Item(31, 3).Value = 6.089783191680908
$sheet0302.Cells.Item(32, 2).Value = 6.215961238881561
$sheet0302.Cells.Item(32, 3).Value = 7.117432594299316
$sheet0302.Cells.Item(33, 2).Value = -0.6436846864172514
$sheet0302.Cells.Item(33, 3).Value = 3.104370832443237
$sheet0302.Cells.Item(34, 2).Value = -0.04236529040875592
$sheet0302.Cells.Item(34, 3).Value = 0.1708746999502182
$sheet0302.Cells.Item(35, 2).Value = -0.7113212780000966
$sheet0302.Cells.Item(35, 3).Value = 2.56526517868042
$sheet0302.Cells.Item(36, 2).Value = -1.575194711333626
$sheet0302.Cells.Item(36, 3).Value = -0.825940728187561

# ---- Restore the originally active sheet/tab ----
$wb.Worksheets.Item("02-28-2024").Activate()

Write-Output "edit complete"
